$d = $word.ActiveDocument

# The document currently has paragraph 1 ("Hallo") and an empty paragraph 2.
# Put "Hola" into the empty paragraph, then add a new paragraph "bonjour" after it.
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Hola"
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs(3)
$p3.Range.Text = "bonjour"
